# Insert a new data row at row 171 (pushing the existing rows 171..251
# down to 172..252) and populate the new row 171 with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(171).Insert()

$ws.Range("A171").Value = 6
$ws.Range("B171").Value = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C171").Value = 'Metropolitana'
$ws.Range("D171").Value = 44839
$ws.Range("E171").Value = 13
$ws.Range("F171").Value = 100112001
$ws.Range("G171").Value = 'Berenjena'
$ws.Range("H171").Value = 'Sin especificar'
$ws.Range("I171").Value = 'Primera'
$ws.Range("J171").Value = 250
$ws.Range("K171").Value = 10000
$ws.Range("L171").Value = 12000
$ws.Range("M171").Value = 10800
$ws.Range("N171").Value = '$/caja 50 unidades'
$ws.Range("O171").Value = 'Región de Arica y Parinacota'
$ws.Range("P171").Value = 216
$ws.Range("Q171").Value = 50
$ws.Range("R171").Value = 'Hortaliza'
